$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("E10").Value = 529
$ws.Range("F10").Value = 258
$ws.Range("H10").Value = 352

# Row 11
$ws.Range("E11").Value = 344

# Row 12
$ws.Range("E12").Value = 520
$ws.Range("F12").Value = 280
$ws.Range("H12").Value = 365

# Row 13
$ws.Range("E13").Value = 133
$ws.Range("F13").Value = 70
$ws.Range("H13").Value = 104

# Row 14
$ws.Range("F14").Value = 66
$ws.Range("G14").Value = 34

# Row 15
$ws.Range("F15").Value = 71
$ws.Range("G15").Value = 50

# Row 22
$ws.Range("E22").Value = 165

# Row 26
$ws.Range("E26").Value = 153

# Row 27
$ws.Range("E27").Value = 322

# Row 30
$ws.Range("E30").Value = 204

# Row 32
$ws.Range("E32").Value = 179
$ws.Range("F32").Value = 106
$ws.Range("H32").Value = 144

# Row 33
$ws.Range("E33").Value = 291
$ws.Range("F33").Value = 146
$ws.Range("H33").Value = 235

# Row 35
$ws.Range("E35").Value = 146

# Row 36
$ws.Range("E36").Value = 70

# Row 38
$ws.Range("E38").Value = 90

# Row 39
$ws.Range("E39").Value = 179

# Row 40
$ws.Range("E40").Value = 257

# Row 42
$ws.Range("F42").Value = 206
$ws.Range("H42").Value = 266

# Row 43
$ws.Range("E43").Value = 115

# Row 45
$ws.Range("E45").Value = 142

# Row 46
$ws.Range("E46").Value = 316
$ws.Range("F46").Value = 171
$ws.Range("H46").Value = 234

# Row 48
$ws.Range("E48").Value = 206
$ws.Range("F48").Value = 89
$ws.Range("H48").Value = 133

# Row 49
$ws.Range("E49").Value = 284

# Row 50
$ws.Range("E50").Value = 237
$ws.Range("G50").Value = 72
$ws.Range("H50").Value = 186

# Row 51
$ws.Range("E51").Value = 230
$ws.Range("F51").Value = 101
$ws.Range("G51").Value = 73
$ws.Range("H51").Value = 174
